$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format so numeric-looking strings
# like "1.005" or "198.10" are preserved exactly (not coerced to numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.956.87"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.641.59"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "215.37"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "0.5081"
$ws.Range("E6").Value = "  +0.76%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.2560"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "0.06381"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "19.48"
$ws.Range("E10").Value = "  -1.20%  "
$ws.Range("D11").Value = "0.07779"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "4.297"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "1.639.79"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "0.5471"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "0.0₅7847"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "64.44"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "26.016.39"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "198.10"
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "9.966"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "6.062"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "1.873"
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "0.1149"
$ws.Range("E26").Value = "  +0.92%  "
$ws.Range("D27").Value = "6.876"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").Value = "15.73"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "1.239"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").Value = "0.05023"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").Value = "3.261"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "2.363"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").Value = "2.592"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("D37").Value = "1.133.42"
$ws.Range("E37").Value = "  -2.54%  "
$ws.Range("D38").Value = "0.5503"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("E39").Value = "  +15.07%  "
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").Value = "1.004"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").Value = "2.542"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").Value = "0.8177"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").Value = "100.19"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "1.779.97"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "0.4530"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "54.91"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "0.05072"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "1.006"
$ws.Range("E51").Value = "  +0.18%  "
